# Scheduled-runner update: refresh market/profit figures (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 197
$ws.Range("I33").Value = 197
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 197
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 32
$ws.Range("N33").ClearContents()

$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H64").Value = 20000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 20000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 20000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 20000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716
$ws.Range("M67").ClearContents()

$ws.Range("H88").Value = 2920.6667
$ws.Range("I88").Value = 2803
$ws.Range("K88").Value = 2803
$ws.Range("M88").Value = -2397

$ws.Range("H91").Value = 2920.6667
$ws.Range("I91").Value = 2803
$ws.Range("K91").Value = 2803
$ws.Range("M91").Value = -1399

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 8547.75

$ws.Range("H137").Value = 1484.9286
$ws.Range("J137").Value = 1399.6666
$ws.Range("L137").Value = 4198.9998
$ws.Range("N137").Value = -9298.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 867.4761999999999
$ws.Range("I97").Value = 599.46155
$ws.Range("K97").Value = 599.46155
$ws.Range("M97").Value = -103.46155

$ws.Range("H122").Value = 1474.1666
$ws.Range("I122").Value = 1420.7778
$ws.Range("J122").Value = 1634.3334
$ws.Range("K122").Value = 4262.3334
$ws.Range("L122").Value = 4903.0002
$ws.Range("M122").Value = -1812.3334
$ws.Range("N122").Value = -9803.0002

$ws.Range("H132").Value = 773.4286
$ws.Range("I132").Value = 819
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 2457
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = 73
$ws.Range("N132").Value = -6560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 19471
$ws.Range("I26").Value = 19471
$ws.Range("K26").Value = 19471
$ws.Range("M26").Value = -19179

$ws.Range("H96").Value = 20000
$ws.Range("I96").Value = 20000
$ws.Range("K96").Value = 20000
$ws.Range("M96").Value = -17254

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 494.2857

$ws.Range("H58").Value = 2000
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 2000
$ws.Range("M58").Value = -1797

$ws.Range("H80").Value = 42606.285
$ws.Range("I80").Value = 45540.668
$ws.Range("K80").Value = 45540.668
$ws.Range("M80").Value = -44417.668

$ws.Range("H83").Value = 42606.285
$ws.Range("I83").Value = 45540.668
$ws.Range("K83").Value = 136622.004
$ws.Range("M83").Value = -131006.004

$ws.Range("H87").Value = 45999
$ws.Range("I87").Value = 45999
$ws.Range("K87").Value = 45999
$ws.Range("M87").Value = -44813

$ws.Range("H90").Value = 45999
$ws.Range("I90").Value = 45999
$ws.Range("K90").Value = 137997
$ws.Range("M90").Value = -132069

$ws.Range("H105").Value = 6174.8335
$ws.Range("I105").Value = 2025
$ws.Range("J105").Value = 8249.75
$ws.Range("K105").Value = 2025
$ws.Range("L105").Value = 8249.75
$ws.Range("M105").Value = -278
$ws.Range("N105").Value = -11743.75

$ws.Range("H122").Value = 1073.2222
$ws.Range("I122").Value = 526.5
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 1579.5
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = 870.5
$ws.Range("N122").Value = -11400.0001

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1494.125
$ws.Range("J5").Value = 1308.3334
$ws.Range("L5").Value = 3925.0002
$ws.Range("N5").Value = -4149.0002

$ws.Range("H58").Value = 500
$ws.Range("J58").Value = 500
$ws.Range("L58").Value = 1500
$ws.Range("N58").Value = -1756

$ws.Range("H69").Value = 1199.1666
$ws.Range("I69").Value = 1199.1666
$ws.Range("K69").Value = 3597.4998
$ws.Range("M69").Value = -2786.4998

$ws.Range("H72").Value = 1199.1666
$ws.Range("I72").Value = 1199.1666
$ws.Range("K72").Value = 10792.4994
$ws.Range("M72").Value = -6736.499400000001

$ws.Range("H113").Value = 543.875
$ws.Range("J113").Value = 699.6667
$ws.Range("L113").Value = 2099.0001
$ws.Range("N113").Value = -6439.0001

$ws.Range("H135").Value = 1494.125
$ws.Range("J135").Value = 1308.3334
$ws.Range("L135").Value = 11775.0006
$ws.Range("N135").Value = -16845.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576

$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996

$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984

$ws.Range("H132").Value = 4290.926
$ws.Range("I132").Value = 4108.095
$ws.Range("K132").Value = 12324.285
$ws.Range("M132").Value = -9794.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3561.8235
$ws.Range("I46").Value = 3533.4443
$ws.Range("J46").Value = 3593.75
$ws.Range("K46").Value = 3533.4443
$ws.Range("L46").Value = 3593.75
$ws.Range("M46").Value = -3345.4443
$ws.Range("N46").Value = -3969.75

$ws.Range("H55").Value = 975.25
$ws.Range("J55").Value = 500
$ws.Range("L55").Value = 500
$ws.Range("N55").Value = -846

$ws.Range("H68").Value = 2499.3333
$ws.Range("I68").Value = 2499.3333
$ws.Range("K68").Value = 2499.3333
$ws.Range("M68").Value = -1750.3333

$ws.Range("H71").Value = 2499.3333
$ws.Range("I71").Value = 2499.3333
$ws.Range("K71").Value = 12496.6665
$ws.Range("M71").Value = -8752.666499999999

$ws.Range("H93").Value = 968.25
$ws.Range("I93").Value = 937.5
$ws.Range("J93").Value = 999
$ws.Range("K93").Value = 937.5
$ws.Range("L93").Value = 999
$ws.Range("M93").Value = 310.5
$ws.Range("N93").Value = -3495

$ws.Range("H120").Value = 100698
$ws.Range("J120").Value = 100698
$ws.Range("L120").Value = 100698
$ws.Range("N120").Value = -110374

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
